$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store numeric-looking values as
# literal text in the source data. Force those specific cells to keep a
# text number format before assigning them, so Excel does not silently
# convert the new values ("244.26", "-0.69%", ...) into real numbers.
$textCells = @(
  "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "E22", "D23", "E23", "E24", "E25", "D26", "E26", "E27", "E28", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "E46", "E47", "D48", "E48", "E49", "E50"
)
foreach ($addr in $textCells) {
  $ws.Range($addr).NumberFormat = "@"
}

# Apply the refreshed symbol list / price / volume values.
$ws.Range("D2").Value = "244.26"
$ws.Range("E2").Value = "-0.69%"
$ws.Range("D3").Value = "27.22"
$ws.Range("E3").Value = "4.15%"
$ws.Range("D4").Value = "5.162"
$ws.Range("E4").Value = "1.27%"
$ws.Range("D5").Value = "0.05636"
$ws.Range("E5").Value = "0.73%"
$ws.Range("D6").Value = "6.480"
$ws.Range("E6").Value = "-0.33%"
$ws.Range("D7").Value = "0.8156"
$ws.Range("E7").Value = "0.59%"
$ws.Range("D8").Value = "0.8336"
$ws.Range("E8").Value = "-0.93%"
$ws.Range("D9").Value = "0.1329"
$ws.Range("E9").Value = "-1.40%"
$ws.Range("D10").Value = "0.06880"
$ws.Range("E10").Value = "-1.80%"
$ws.Range("D11").Value = "0.02936"
$ws.Range("E11").Value = "3.64%"
$ws.Range("D12").Value = "0.09391"
$ws.Range("E12").Value = "0.03%"
$ws.Range("D13").Value = "0.001508"
$ws.Range("E13").Value = "-0.90%"
$ws.Range("D14").Value = "0.04246"
$ws.Range("E14").Value = "-9.51%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "0.0005967"
$ws.Range("E15").Value = "-0.54%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.006118"
$ws.Range("E16").Value = "-0.32%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.534"
$ws.Range("E17").Value = "-0.61%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "3.007"
$ws.Range("E18").Value = "-0.55%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "2.227"
$ws.Range("E19").Value = "5.12%"
$ws.Range("E22").Value = "-2.16%"
$ws.Range("D23").Value = "3.751"
$ws.Range("E23").Value = "-0.07%"
$ws.Range("E24").Value = "-0.11%"
$ws.Range("E25").Value = "-1.90%"
$ws.Range("D26").Value = "0.004473"
$ws.Range("E26").Value = "-3.04%"
$ws.Range("E27").Value = "2.07%"
$ws.Range("E28").Value = "-0.48%"
$ws.Range("D40").Value = "0.03652"
$ws.Range("E40").Value = "-0.08%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "0.1365"
$ws.Range("E41").Value = "29.58%"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "0.006048"
$ws.Range("E42").Value = "-1.16%"
$ws.Range("D43").Value = "0.002629"
$ws.Range("E43").Value = "5.18%"
$ws.Range("D44").Value = "0.008233"
$ws.Range("E44").Value = "-5.43%"
$ws.Range("D45").Value = "0.00005383"
$ws.Range("E45").Value = "1.69%"
$ws.Range("E46").Value = "-0.02%"
$ws.Range("E47").Value = "1.62%"
$ws.Range("D48").Value = "0.002649"
$ws.Range("E48").Value = "29.10%"
$ws.Range("E49").Value = "-0.02%"
$ws.Range("E50").Value = "-0.02%"
